$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.805.19'
$ws.Range("E2").Value = '  +0.93%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.954.83'
$ws.Range("E3").Value = '  +3.77%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.54'
$ws.Range("E5").Value = '  +3.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5987'
$ws.Range("E6").Value = '  +27.83%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9994'
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3144'
$ws.Range("E8").Value = '  +8.70%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '24.50'
$ws.Range("E9").Value = '  +9.93%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06870'
$ws.Range("E10").Value = '  +6.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8060'
$ws.Range("E11").Value = '  +10.55%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '101.34'
$ws.Range("E12").Value = '  +6.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07967'
$ws.Range("E13").Value = '  +2.90%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.927.86'
$ws.Range("E14").Value = '  +2.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.350'
$ws.Range("E15").Value = '  +3.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '281.12'
$ws.Range("E16").Value = '  -0.49%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.790.79'
$ws.Range("E17").Value = '  +0.91%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.74'
$ws.Range("E18").Value = '  +5.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007696'
$ws.Range("E19").Value = '  +2.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.596'
$ws.Range("E20").Value = '  +6.60%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.177.86'
$ws.Range("E21").Value = '  +2.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9997'
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9990'
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.640'
$ws.Range("E24").Value = '  +6.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.469'
$ws.Range("E25").Value = '  +4.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.26'
$ws.Range("E26").Value = '  +1.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.55'
$ws.Range("E27").Value = '  +3.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.077'
$ws.Range("E28").Value = '  +9.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1111'
$ws.Range("E29").Value = '  +14.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.353'
$ws.Range("E30").Value = '  +1.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.544'
$ws.Range("E31").Value = '  +4.84%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.474'
$ws.Range("E32").Value = '  +4.70%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.345'
$ws.Range("E33").Value = '  +5.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04989'
$ws.Range("E34").Value = '  +2.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.176'
$ws.Range("E35").Value = '  +4.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7221'
$ws.Range("E36").Value = '  +4.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.715'
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01963'
$ws.Range("E38").Value = '  +3.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.924'
$ws.Range("E39").Value = '  +3.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '77.67'
$ws.Range("E40").Value = '  +2.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.415'
$ws.Range("E41").Value = '  +4.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4519'
$ws.Range("E42").Value = '  +6.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.020'
$ws.Range("E43").Value = '  +0.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8442'
$ws.Range("E44").Value = '  +2.70%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.73'
$ws.Range("E46").Value = '  +1.44%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.04'
$ws.Range("E47").Value = '  +5.57%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.282'
$ws.Range("E48").Value = '  +4.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.97'
$ws.Range("E49").Value = '  +2.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4142'
$ws.Range("E50").Value = '  +5.32%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '917.56'
$ws.Range("E51").Value = '  +0.53%  '
